# US56, Task75 - Session Handling for admin portal
# Append three new user rows (two "test" accounts + one "admin" account)
# to the "Users" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - test user
$ws.Range("A3").Value = "test"
$ws.Range("B3").Value = "name"
$ws.Range("C3").Value = "abc@gmail.com"
$ws.Range("D3").Value = "student"
$ws.Range("E3").Value = "high"
$ws.Range("F3").Value = "pwd*"

# Row 4 - test user (duplicate entry)
$ws.Range("A4").Value = "test"
$ws.Range("B4").Value = "name"
$ws.Range("C4").Value = "abc@gmail.com"
$ws.Range("D4").Value = "student"
$ws.Range("E4").Value = "high"
$ws.Range("F4").Value = "pwd*"

# Row 5 - admin user (Grade column intentionally blank/empty string)
$ws.Range("A5").Value = "admin"
$ws.Range("B5").Value = "admin"
$ws.Range("C5").Value = "admin@gmail.com"
$ws.Range("D5").Value = "admin"
$ws.Range("E5").Value = "'"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "pwd*"
